# Fix Training Data Issue (#48)
# The "Date" column (BF) held strings like "5-31-2013-14" which were
# actually off by a day because of how the source NBA stats were
# reported. Correct each row's date text to ISO form "2014-05-31".
#
# The column is pre-formatted as Text so Excel stores the corrected
# value as the literal string "2014-05-31" instead of silently
# re-interpreting it as a date serial number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 31
$dateCol = 58   # column BF

$rng = $ws.Range($ws.Cells.Item($firstRow, $dateCol), $ws.Cells.Item($lastRow, $dateCol))
$rng.NumberFormat = "@"

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, $dateCol).Value = "2014-05-31"
}
